$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 currently holds the oldest record (date 2021-12-27 / serial 44557).
# A new weekly record is appended, so that old record is pushed down to a
# new row 21 with its original values intact, and row 20 is overwritten
# with the new weekly data point.

# New row 21: copy of the previous row 20 contents.
$ws.Range("A21").Value = 8
$ws.Range("B21").Value = "Terminal La Palmera de La Serena"
$ws.Range("C21").Value = "Coquimbo"
$ws.Range("D21").Value = 44557
$ws.Range("D21").NumberFormat = $ws.Range("D20").NumberFormat
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = 100114002
$ws.Range("G21").Value = "Camote"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 400
$ws.Range("K21").Value = 13000
$ws.Range("L21").Value = 14000
$ws.Range("M21").Value = 13500
$ws.Range("N21").Value = "$/malla 18 kilos"
$ws.Range("O21").Value = "Perú"
$ws.Range("P21").Value = 750
$ws.Range("Q21").Value = 18
$ws.Range("R21").Value = "Hortaliza"

# Overwrite row 20 with the new weekly data point (date 2023-09-25).
$ws.Range("D20").Value = 45194
$ws.Range("K20").Value = 16500
$ws.Range("L20").Value = 17000
$ws.Range("M20").Value = 16750
$ws.Range("P20").Value = 931
